$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.025.77"
$ws.Range("E2").Value = "  -0.23%  "
$ws.Range("D3").Value = "2.416.27"
$ws.Range("E3").Value = "  -0.74%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "552.99"
$ws.Range("E5").Value = "  -0.66%  "
$ws.Range("D6").Value = "137.30"
$ws.Range("E6").Value = "  -1.23%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  +3.89%  "
$ws.Range("E9").Value = "  -1.94%  "
$ws.Range("E10").Value = "  -2.42%  "
$ws.Range("D11").Value = "0.147"
$ws.Range("E11").Value = "  -0.90%  "
$ws.Range("E12").Value = "  -2.17%  "
$ws.Range("D13").Value = "25.42"
$ws.Range("E13").Value = "  +2.11%  "
$ws.Range("D14").Value = "2.847.96"
$ws.Range("E14").Value = "  -0.67%  "
$ws.Range("D15").Value = "59.947.52"
$ws.Range("E15").Value = "  -0.14%  "
$ws.Range("E16").Value = "  -1.95%  "
$ws.Range("D17").Value = "2.411.27"
$ws.Range("E17").Value = "  -0.81%  "
$ws.Range("D18").Value = "11.32"
$ws.Range("E18").Value = "  -1.39%  "
$ws.Range("D19").Value = "4.43"
$ws.Range("E19").Value = "  -0.40%  "
$ws.Range("D20").Value = "329.03"
$ws.Range("E20").Value = "  -1.90%  "
$ws.Range("D21").Value = "6.67"
$ws.Range("E21").Value = "  -3.78%  "
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("D23").Value = "65.93"
$ws.Range("E23").Value = "  +1.79%  "
$ws.Range("D24").Value = "0.176"
$ws.Range("E24").Value = "  +3.34%  "
$ws.Range("D25").Value = "8.63"
$ws.Range("E25").Value = "  +0.85%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("E27").Value = "  -0.08%  "
$ws.Range("D28").Value = "0.0₃0777"
$ws.Range("E28").Value = "  -2.02%  "
$ws.Range("E29").Value = "  -2.39%  "
$ws.Range("D30").Value = "168.99"
$ws.Range("E30").Value = "  -1.37%  "
$ws.Range("D31").Value = "6.06"
$ws.Range("E31").Value = "  -4.06%  "
$ws.Range("E32").Value = "  -1.02%  "
$ws.Range("E33").Value = "  -0.33%  "
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("D37").Value = "4.19"
$ws.Range("E37").Value = "  -1.89%  "
$ws.Range("E38").Value = "  -2.55%  "
$ws.Range("D39").Value = "322.21"
$ws.Range("E39").Value = "  +1.28%  "
$ws.Range("E40").Value = "  -2.90%  "
$ws.Range("E41").Value = "  -1.79%  "
$ws.Range("D42").Value = "140.42"
$ws.Range("E42").Value = "  -1.99%  "
$ws.Range("E43").Value = "  +0.79%  "
$ws.Range("D44").Value = "19.61"
$ws.Range("E44").Value = "  +0.86%  "
$ws.Range("E45").Value = "  -1.88%  "
$ws.Range("D46").Value = "0.580"
$ws.Range("E46").Value = "  +1.15%  "
$ws.Range("E47").Value = "  -1.72%  "
$ws.Range("E48").Value = "  -6.19%  "
$ws.Range("E50").Value = "  -3.51%  "
$ws.Range("E51").Value = "  -1.04%  "
